# TC14 Trials Filter Diagnosis-Hurthle: add a "StatQuery" column
# Inserts a new column B (header "StatQuery") holding a second Cypher
# query that returns file/case/trial counts, shifting the former
# "dbExcel" / "WebExcel" columns from B/C to C/D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing column B ("dbExcel"),
# pushing dbExcel -> C and WebExcel -> D.
$ws.Columns("B").Insert()

# Match column A's width for the newly inserted column B.
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# New header for the inserted column.
$ws.Range("B1").Value = "StatQuery"

# New query text for row 2 (wrap-text style carries over automatically
# from the Insert, matching A2's style).
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Hurthle cell neoplasm (thyroid)'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Make sure the wrap-text style used for A2 is applied to B2 as well.
$ws.Range("B2").WrapText = $true

# Move the active selection to A4, matching the saved view state.
$ws.Range("A4").Select()
